$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 19228.334
$ws.Range("J69").Value = 19635
$ws.Range("L69").Value = 58905
$ws.Range("N69").Value = -60653
$ws.Range("H70").Value = 125002160
$ws.Range("J70").Value = 142859420
$ws.Range("L70").Value = 428578260
$ws.Range("N70").Value = -428578800
$ws.Range("H72").Value = 19228.334
$ws.Range("J72").Value = 19635
$ws.Range("L72").Value = 176715
$ws.Range("N72").Value = -185451
$ws.Range("H73").Value = 125002160
$ws.Range("J73").Value = 142859420
$ws.Range("L73").Value = 428578260
$ws.Range("N73").Value = -428580132
$ws.Range("H86").Value = 3654.6155
$ws.Range("I86").Value = 3457.889
$ws.Range("K86").Value = 3457.889
$ws.Range("M86").Value = -2334.889
$ws.Range("H89").Value = 3654.6155
$ws.Range("I89").Value = 3457.889
$ws.Range("K89").Value = 17289.445
$ws.Range("M89").Value = -11673.445
$ws.Range("H106").Value = 5000.885
$ws.Range("I106").Value = 5332.864
$ws.Range("K106").Value = 5332.864
$ws.Range("M106").Value = -4701.864
$ws.Range("H112").Value = 3929.6316
$ws.Range("J112").Value = 3929.6316
$ws.Range("L112").Value = 11788.8948
$ws.Range("N112").Value = -14004.8948
$ws.Range("H129").Value = 1797.6
$ws.Range("I129").Value = 1215.875
$ws.Range("K129").Value = 3647.625
$ws.Range("M129").Value = 1352.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9389.277
$ws.Range("I32").Value = 1971.9286
$ws.Range("K32").Value = 1971.9286
$ws.Range("M32").Value = -1684.9286
$ws.Range("H45").Value = 2940.56
$ws.Range("J45").Value = 3316.3333
$ws.Range("L45").Value = 3316.3333
$ws.Range("N45").Value = -4070.3333
$ws.Range("H61").Value = 12714.913
$ws.Range("I61").Value = 8896.714
$ws.Range("J61").Value = 18654.334
$ws.Range("K61").Value = 8896.714
$ws.Range("L61").Value = 18654.334
$ws.Range("M61").Value = -8684.714
$ws.Range("N61").Value = -19078.334
$ws.Range("H74").Value = 4484.4287
$ws.Range("I74").Value = 3798
$ws.Range("K74").Value = 3798
$ws.Range("M74").Value = -2924
$ws.Range("H77").Value = 4484.4287
$ws.Range("I77").Value = 3798
$ws.Range("K77").Value = 18990
$ws.Range("M77").Value = -14622
$ws.Range("H132").Value = 1934.6923
$ws.Range("I132").Value = 1653.7542
$ws.Range("K132").Value = 4961.2626
$ws.Range("M132").Value = -2431.2626
$ws.Range("H136").Value = 12714.913
$ws.Range("I136").Value = 8896.714
$ws.Range("J136").Value = 18654.334
$ws.Range("K136").Value = 26690.142
$ws.Range("L136").Value = 55963.00199999999
$ws.Range("M136").Value = -24140.142
$ws.Range("N136").Value = -61063.00199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4169572
$ws.Range("I86").Value = 9261356
$ws.Range("J86").Value = 3566.818
$ws.Range("K86").Value = 9261356
$ws.Range("L86").Value = 3566.818
$ws.Range("M86").Value = -9260233
$ws.Range("N86").Value = -5812.818
$ws.Range("H89").Value = 4169572
$ws.Range("I89").Value = 9261356
$ws.Range("J89").Value = 3566.818
$ws.Range("K89").Value = 46306780
$ws.Range("L89").Value = 17834.09
$ws.Range("M89").Value = -46301164
$ws.Range("N89").Value = -29066.09

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4752.6562
$ws.Range("I31").Value = 4467.345
$ws.Range("K31").Value = 4467.345
$ws.Range("M31").Value = -4172.345
$ws.Range("H34").Value = 4752.6562
$ws.Range("I34").Value = 4467.345
$ws.Range("K34").Value = 4467.345
$ws.Range("M34").Value = -4265.345
$ws.Range("H58").Value = 2041.9615
$ws.Range("I58").Value = 2155.348
$ws.Range("J58").Value = 1172.6666
$ws.Range("K58").Value = 2155.348
$ws.Range("L58").Value = 1172.6666
$ws.Range("M58").Value = -1952.348
$ws.Range("N58").Value = -1578.6666
$ws.Range("H107").Value = 333.26315
$ws.Range("I107").Value = 302.4375
$ws.Range("K107").Value = 302.4375
$ws.Range("M107").Value = 1617.5625
$ws.Range("H132").Value = 1371.3572
$ws.Range("I132").Value = 1359.1621
$ws.Range("K132").Value = 4077.4863
$ws.Range("M132").Value = -1547.4863
$ws.Range("H134").Value = 1113.1136
$ws.Range("I134").Value = 963.7619
$ws.Range("K134").Value = 2891.2857
$ws.Range("M134").Value = -356.2856999999999
$ws.Range("H136").Value = 2041.9615
$ws.Range("I136").Value = 2155.348
$ws.Range("J136").Value = 1172.6666
$ws.Range("K136").Value = 6466.044
$ws.Range("L136").Value = 3517.9998
$ws.Range("M136").Value = -3916.044
$ws.Range("N136").Value = -8617.9998
$ws.Range("H137").Value = 53144.5
$ws.Range("I137").Value = 48262.5
$ws.Range("J137").Value = 54771.832
$ws.Range("K137").Value = 48262.5
$ws.Range("L137").Value = 54771.832
$ws.Range("M137").Value = -43162.5
$ws.Range("N137").Value = -64971.832

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 938.5333000000001
$ws.Range("I5").Value = 312.7143
$ws.Range("K5").Value = 938.1428999999999
$ws.Range("M5").Value = -826.1428999999999
$ws.Range("H43").Value = 17500
$ws.Range("J43").Value = 17500
$ws.Range("L43").Value = 52500
$ws.Range("N43").Value = -52728
$ws.Range("H45").Value = 999.5
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 999.5
$ws.Range("K45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("M45").Value = 2998.5
$ws.Range("N45").Value = -4062.5
$ws.Range("H63").Value = 499
$ws.Range("I63").Value = 499
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1497
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -748
$ws.Range("H66").Value = 499
$ws.Range("I66").Value = 499
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 4491
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -747
$ws.Range("H116").Value = 1727.5
$ws.Range("I116").Value = 1727.5
$ws.Range("K116").Value = 5182.5
$ws.Range("M116").Value = -1740.5
$ws.Range("H135").Value = 938.5333000000001
$ws.Range("I135").Value = 312.7143
$ws.Range("K135").Value = 2814.4287
$ws.Range("M135").Value = -279.4286999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5324.9546
$ws.Range("I46").Value = 6171.2856
$ws.Range("J46").Value = 4930
$ws.Range("K46").Value = 6171.2856
$ws.Range("L46").Value = 4930
$ws.Range("M46").Value = -5983.2856
$ws.Range("N46").Value = -5306
$ws.Range("H54").Value = 70000
$ws.Range("J54").Value = 70000
$ws.Range("L54").Value = 70000
$ws.Range("N54").Value = -71288
$ws.Range("H61").Value = 991.3333
$ws.Range("I61").Value = 991.3333
$ws.Range("K61").Value = 991.3333
$ws.Range("M61").Value = -789.3333
$ws.Range("H113").Value = 991.3333
$ws.Range("I113").Value = 991.3333
$ws.Range("K113").Value = 991.3333
$ws.Range("M113").Value = 1178.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3369666.8
$ws.Range("I3").Value = 3369666.8
$ws.Range("K3").Value = 3369666.8
$ws.Range("M3").Value = -3369552.8
$ws.Range("H58").Value = 17996
$ws.Range("I58").Value = 17999
$ws.Range("J58").Value = 17990
$ws.Range("K58").Value = 17999
$ws.Range("L58").Value = 17990
$ws.Range("M58").Value = -17691
$ws.Range("N58").Value = -18606
$ws.Range("H132").Value = 1875.3112
$ws.Range("I132").Value = 1903.3903
$ws.Range("K132").Value = 5710.1709
$ws.Range("M132").Value = -3180.1709
